$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.324.72'
$ws.Range("E2").Value = '  +0.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.082.47'
$ws.Range("E3").Value = '  +4.00%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.86'
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.31'
$ws.Range("E6").Value = '  +3.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.080.52'
$ws.Range("E8").Value = '  +4.08%  '

$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  +0.56%  '

$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.479'
$ws.Range("E12").Value = '  +5.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.65'
$ws.Range("E14").Value = '  +6.94%  '

$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.591.81'
$ws.Range("E16").Value = '  +3.92%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.361.13'
$ws.Range("E17").Value = '  +0.35%  '

$ws.Range("E18").Value = '  +4.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.082.94'
$ws.Range("E19").Value = '  +4.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.08'
$ws.Range("E20").Value = '  +17.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '463.49'
$ws.Range("E21").Value = '  +3.45%  '

$ws.Range("E22").Value = '  +5.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.44'
$ws.Range("E23").Value = '  +3.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.15'
$ws.Range("E24").Value = '  +1.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.81'
$ws.Range("E25").Value = '  +5.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("E26").Value = '  +2.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  +1.09%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  +1.06%  '

$ws.Range("E31").Value = '  +3.18%  '

$ws.Range("E32").Value = '  +2.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.33'
$ws.Range("E33").Value = '  +4.21%  '

$ws.Range("E34").Value = '  +5.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +2.12%  '

$ws.Range("E37").Value = '  +3.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.67'
$ws.Range("E38").Value = '  +12.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.10'
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("E40").Value = '  +4.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.04'
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.122'
$ws.Range("E42").Value = '  +2.53%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("E43").Value = '  +2.08%  '

$ws.Range("E44").Value = '  +3.51%  '

$ws.Range("E45").Value = '  +1.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '384.03'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.774.54'
$ws.Range("E47").Value = '  +2.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.91'
$ws.Range("E48").Value = '  +3.12%  '

$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("E50").Value = '  +6.09%  '

$ws.Range("E51").Value = '  +4.67%  '
